$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B=0.1288888888888889; C=31.10111111111111; D=5.977777777777778; E=84.48222222222222; F=50.21555555555555; G=97.5911111111111; H=91.98444444444445; I=57.41333333333333; J=41.22111111111111; K=19.66888888888889; L=5.138888888888889; M=48.74333333333334; N=15.10666666666667; O=56.51888888888889; P=96.52777777777777 }
    3 = @{ B=0.08; C=35.30714285714286; D=6.351428571428572; E=83.30285714285715; F=58.17000000000001; G=95.34428571428573; H=88.96571428571428; I=52.73285714285714; J=45.15857142857143; K=18.90142857142857; L=6.300000000000002; M=57.08714285714285; N=13.50571428571429; O=56.52142857142857; P=96.90428571428572 }
    4 = @{ B=0.5675; C=42.66; D=8.43; E=86.91249999999999; F=54.8; G=94.06999999999999; H=91.55000000000001; I=53.785; J=42.8725; K=23.76; L=10.9925; M=57.63; N=17.6875; O=54.03250000000001; P=97.77249999999999 }
    5 = @{ B=0.2625; C=41.145; D=4.8425; E=80.205; F=46.8475; G=97.3725; H=86.37; I=59.415; J=35.5925; K=14.09; L=5.095; M=54.8625; N=10.7375; O=59.39; P=97.345 }
    6 = @{ B=0.1233333333333333; C=43.09; D=12.59333333333333; E=86.85666666666667; F=53.40666666666667; G=95.46666666666665; H=88.25999999999999; I=57.97666666666667; J=42.83666666666667; K=15.79333333333333; L=3.9; M=52.62; N=12.63666666666667; O=57.27; P=98.78666666666668 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
